# #5: property aircraft done
# The "property_category" column was hard-coded to "land" when the report
# was generated, regardless of which asset-type sheet it lived on. This
# fixes the mislabeled rows on the "building" (建物) and "car" (汽車) sheets
# so their property_category matches the sheet they actually belong to.

$wb = $excel.ActiveWorkbook

# Sheet 2 = "建物" (building): column I = property_category, data rows 2-9
$wsBuilding = $wb.Worksheets.Item(2)
for ($row = 2; $row -le 9; $row++) {
    $wsBuilding.Cells.Item($row, 9).Value = "building"
}

# Sheet 3 = "汽車" (car): column H = property_category, data rows 2-3
$wsCar = $wb.Worksheets.Item(3)
for ($row = 2; $row -le 3; $row++) {
    $wsCar.Cells.Item($row, 8).Value = "car"
}
